$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H2: report date text string
$ws.Range("H2").Value = "2020-09-30 00:00:00"

# I2: basic eps
$ws.Range("I2").Value = 0.54

# K2: total operate income
$ws.Range("K2").Value = 99369292.31999999

# L2: parent netprofit
$ws.Range("L2").Value = 16094268.55

# N2: YSTZ (was empty, now numeric)
$ws.Range("N2").Value = 122.9209529084

# O2: SJLTZ (was empty, now numeric)
$ws.Range("O2").Value = 348.5532237172

# P2: BPS (was empty, now numeric)
$ws.Range("P2").Value = 3.307587800667

# Q2: MGJYXJJE (was empty, now numeric)
$ws.Range("Q2").Value = -0.393075698667

# R2: XSMLL
$ws.Range("R2").Value = 30.2504255874

# AB2: ISNEW flag, keep as text "1" (leading apostrophe forces text entry, like the original)
$ws.Range("AB2").Value = "'1"

# AC2: QDATE
$ws.Range("AC2").Value = "2020Q3"

# AD2: DATATYPE
$ws.Range("AD2").Value = "2020年 三季报"

# AE2: DATAYEAR, keep as text "2020" (leading apostrophe forces text entry, like the original)
$ws.Range("AE2").Value = "'2020"
